$d = $word.ActiveDocument

# Locate the "056/25 ..." entry paragraph by its content (rather than a
# hardcoded index) and remember its 1-based position.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*056/25*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)
# Cloning via InsertParagraphAfter picks up the source paragraph's
# formatting (spacing-after + run size), matching the sibling entries.
$target.Range.InsertParagraphAfter()

# Re-fetch the freshly created paragraph (immediately after the target)
# and fill in its text.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "125/25 – 30.06.2025. – Obuka – 5.000 RSD (kandidat van evidencije)"

# Update the total amount to reflect the added entry.
$d.Content.Find.Execute("Ukupno: 5.000 RSD", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Ukupno: 10.000 RSD", 2)
